# Append two new dev-log entries to the "Logs" sheet, matching the
# existing table's date/wrap-text formatting by copying row 42's format
# down into the two new rows before writing the new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A42:B42").Copy()
$ws.Range("A43:B43").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A44:B44").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A43").Value = 45539
$ws.Range("B43").Value = "add mass spawn of enemies on attack or agr zone enter, add primitive player respawn zone and dead souls drop and life minus"

$ws.Range("A44").Value = 45540
$ws.Range("B44").Value = "add gui theme, add custom gui bars, connect them to player, fix spawner"

# Matches the author's new selection after entering the two rows.
$ws.Range("B46").Select()
